$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the (empty / unused) drawing object that anchors to this sheet,
# mirroring the worksheet no longer referencing xl/drawings/drawing1.xml.
if ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.SelectAll()
    $excel.Selection.Delete()
}

# --- New metadata row (row 2) -------------------------------------------
$ws.Range("A2").Value = "MCH198"
$ws.Range("C2").Value = "PROCEEDINGS OF THE WORLD CONFERENCE AGAINST APARTHEID, RACISM AND COLONIALISM IN SOUTHERN AFRICA, LISBON 16-19 JUNE 1977"

# date_s needs to land as text ("1977"), not a number - build it via a
# formula then freeze it back to a plain value so no quote-prefix / custom
# number-format ends up tagging the cell.
$ws.Range("D2").Formula = '="1977"'
$ws.Range("D2").Value = $ws.Range("D2").Value

$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# --- Formatting for the new row ------------------------------------------
# New font: Calibri 10pt, theme text color - used by A2:H2 (F2 additionally
# keeps an alignment flag).
$dataRow = $ws.Range("A2:H2")
$dataRow.Font.ThemeColor = 1
$dataRow.Font.Name = "Calibri"
$dataRow.Font.Size = 10
